$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.745.50"
$ws.Range("E2").Value = "  +5.60%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.705.63"
$ws.Range("E3").Value = "  +3.24%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9987"
$ws.Range("E4").Value = "  -0.63%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "331.05"
$ws.Range("E5").Value = "  +6.38%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9973"
$ws.Range("E6").Value = "  -0.44%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3686"
$ws.Range("E7").Value = "  +1.21%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "48.52"
$ws.Range("E8").Value = "  +3.13%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3307"
$ws.Range("E9").Value = "  +1.80%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.167"
$ws.Range("E10").Value = "  +3.98%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07349"
$ws.Range("E11").Value = "  +4.53%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9964"
$ws.Range("E12").Value = "  -0.55%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.199"
$ws.Range("E13").Value = "  +4.43%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.97"
$ws.Range("E14").Value = "  +3.11%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.863"
$ws.Range("E15").Value = "  +4.14%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.700.93"
$ws.Range("E16").Value = "  +3.02%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001073"
$ws.Range("E17").Value = "  +2.87%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06620"
$ws.Range("E18").Value = "  +0.19%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "81.05"
$ws.Range("E19").Value = "  +3.95%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9972"
$ws.Range("E20").Value = "  -0.41%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.18"
$ws.Range("E21").Value = "  +3.73%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.045"
$ws.Range("E22").Value = "  +2.02%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.00"
$ws.Range("E23").Value = "  +4.66%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "25.694.58"
$ws.Range("E24").Value = "  +5.46%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.450"
$ws.Range("E25").Value = "  -1.11%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.486"
$ws.Range("E26").Value = "  +6.14%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "149.42"
$ws.Range("E27").Value = "  +1.48%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.19"
$ws.Range("E28").Value = "  +3.40%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.271"
$ws.Range("E29").Value = "  +8.57%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.889.46"

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "128.14"
$ws.Range("E31").Value = "  +3.49%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.094"
$ws.Range("E32").Value = "  +0.66%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.965"
$ws.Range("E33").Value = "  +5.95%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08480"
$ws.Range("E34").Value = "  +0.29%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.683"
$ws.Range("E35").Value = "  +0.97%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.74"
$ws.Range("E36").Value = "  +3.86%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.317"
$ws.Range("E37").Value = "  +2.68%  "

# Row 38
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06236"
$ws.Range("E38").Value = "  +3.50%  "

# Row 39
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.270"
$ws.Range("E39").Value = "  +1.57%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.546"
$ws.Range("E40").Value = "  +5.12%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2122"
$ws.Range("E41").Value = "  +2.88%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.02258"
$ws.Range("E42").Value = "  +1.96%  "

# Row 43
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "14.33"
$ws.Range("E43").Value = "  +13.82%  "

# Row 44
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6108"
$ws.Range("E44").Value = "  +3.87%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9970"
$ws.Range("E45").Value = "  -0.44%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.844"
$ws.Range("E46").Value = "  +2.00%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5854"
$ws.Range("E47").Value = "  +4.43%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.88"
$ws.Range("E48").Value = "  +2.81%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.005"
$ws.Range("E49").Value = "  +3.16%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07219"
$ws.Range("E50").Value = "  +4.90%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.209"
$ws.Range("E51").Value = "  +3.33%  "
